{"js": "// Office.js (Word JavaScript API) script\n// 1) Split \"The red dot is the Hero...\" into two runs: \"The red square\" and\n//    \" is the Hero...\", changing \"dot\" -> \"square\" in the process.\n// 2) Add noProof (w:noProof) to the run that wraps the second picture\n//    (Picture 3) in the document.\n\nconst body = context.document.body;\n\n// --- Change 1: \"The red dot\" -> \"The red square\", split into two runs ---\nconst results = body.search(\"The red dot\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n  hit.insertText(\"The red square\", \"Replace\");\n  // Toggle a character property on just this sub-range so the engine keeps\n  // it as its own run instead of re-merging it with the (identically\n  // formatted) text that follows.\n  hit.font.bold = true;\n  await context.sync();\n  hit.font.bold = false;\n  await context.sync();\n}\n\n// --- Change 2: add <w:noProof/> to the run wrapping the second picture ---\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length >= 2) {\n  const secondPicRange = pics.items[1].getRange();\n  secondPicRange.hasNoProofing = true;\n  await context.sync();\n}\n", "ps1": "# Word COM interop script\n# 1) Split \"The red dot is the Hero...\" into two runs: \"The red square\" and\n#    \" is the Hero...\", changing \"dot\" -> \"square\" in the process.\n# 2) Add noProof (w:noProof) to the run that wraps the second picture\n#    (Picture 3) in the document.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"The red dot\" -> \"The red square\", split into two runs ---\n$hit = $d.Content\n$find = $hit.Find\n$find.ClearFormatting()\n$find.Text = \"The red dot\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 0\n$found = $find.Execute()\n\nif ($found) {\n    $hit.Text = \"The red square\"\n    # Toggle a character property on just this sub-range so Word keeps it as\n    # its own run instead of re-merging it with the (identically formatted)\n    # text that follows.\n    $hit.Font.Bold = $true\n    $hit.Font.Bold = $false\n}\n\n# --- Change 2: add <w:noProof/> to the run wrapping the second picture ---\n$shapes = $d.InlineShapes\nif ($shapes.Count -ge 2) {\n    $pic = $shapes.Item(2)\n    $pic.Range.NoProofing = 1\n}\n"}
